$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# PART 1: Fill the empty paragraph after "...password and ESSID are correct."
#         and insert the new "Logging into the RPi via SSH" section.
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("password and ESSID are correct.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $rng.Paragraphs(1).Next()
$cursor = $targetPara.Range.Start

$rng = $d.Range($cursor, $cursor)
$rng.InsertAfter("If your RPi has connected to your WiFi, you can now move to controlling the RPi via SSH. ")
$rng.Font.Size = 12
$cursor = $rng.End

# New empty paragraph
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1

# New paragraph: "Logging into the RPi via SSH" (bold)
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1
$rng = $d.Range($cursor, $cursor)
$rng.InsertAfter("Logging into the RPi via SSH")
$rng.Font.Size = 12
$rng.Font.Bold = $true
$cursor = $rng.End

# New paragraph: "On your RPi, type ifconfig. ..."
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1
$rng = $d.Range($cursor, $cursor)
$rng.InsertAfter("On your RPi, type ifconfig. Under wlan0, the inet addr is the address of the RPi on your network. ")
$rng.Font.Size = 12
$cursor = $rng.End

# New paragraph: "This will be the address you use to login via ssh. ..."
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1
$rng = $d.Range($cursor, $cursor)
$rng.InsertAfter("This will be the address you use to login via ssh. This could also be obtained via a network scan, but I think this is easier. The login username is pi, and the password is raspberry. ")
$rng.Font.Size = 12
$cursor = $rng.End

# New paragraph: "I recommend the Bitvise SSH client ..."
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1
$rng = $d.Range($cursor, $cursor)
$rng.InsertAfter("I recommend the Bitvise SSH client because it has a built-in SFTP file manager, so you can transfer files easily between the RPi and your desktop. You could also just use Git. ")
$rng.Font.Size = 12
$cursor = $rng.End

# New empty paragraph
$rng = $d.Range($cursor, $cursor)
$rng.InsertParagraphAfter()
$cursor = $cursor + 1

Write-Host "Part 1 done. cursor=$cursor"
